# Auto-generated script applying numeric corrections to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 730.63635
$ws.Range("J98").Value = 1081.1111
$ws.Range("L98").Value = 1081.1111
$ws.Range("N98").Value = -4077.1111
$ws.Range("H122").Value = 730.63635
$ws.Range("J122").Value = 1081.1111
$ws.Range("L122").Value = 3243.3333
$ws.Range("N122").Value = -8143.3333
$ws.Range("H137").Value = 1757.1538
$ws.Range("I137").Value = 2015.6666
$ws.Range("K137").Value = 6046.9998
$ws.Range("M137").Value = -3496.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5150.507
$ws.Range("I32").Value = 4385.9546
$ws.Range("K32").Value = 4385.9546
$ws.Range("M32").Value = -4098.9546
$ws.Range("H61").Value = 1597.3889
$ws.Range("I61").Value = 1549.3125
$ws.Range("J61").Value = 1982
$ws.Range("K61").Value = 1549.3125
$ws.Range("L61").Value = 1982
$ws.Range("M61").Value = -1337.3125
$ws.Range("N61").Value = -2406
$ws.Range("H110").Value = 586.25
$ws.Range("J110").Value = 783
$ws.Range("L110").Value = 783
$ws.Range("N110").Value = -4873
$ws.Range("H122").Value = 1230.0385
$ws.Range("I122").Value = 1070.75
$ws.Range("K122").Value = 3212.25
$ws.Range("M122").Value = -762.25
$ws.Range("H132").Value = 20080.643
$ws.Range("I132").Value = 2002.091
$ws.Range("J132").Value = 86368.664
$ws.Range("K132").Value = 6006.272999999999
$ws.Range("L132").Value = 259105.992
$ws.Range("M132").Value = -3476.272999999999
$ws.Range("N132").Value = -264165.992
$ws.Range("H136").Value = 1597.3889
$ws.Range("I136").Value = 1549.3125
$ws.Range("J136").Value = 1982
$ws.Range("K136").Value = 4647.9375
$ws.Range("L136").Value = 5946
$ws.Range("M136").Value = -2097.9375
$ws.Range("N136").Value = -11046
$ws.Range("H139").Value = 40326.75
$ws.Range("J139").Value = 40326.75
$ws.Range("L139").Value = 40326.75
$ws.Range("N139").Value = -50606.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H134").Value = 2448.8
$ws.Range("I134").Value = 2380.509
$ws.Range("J134").Value = 3200
$ws.Range("K134").Value = 7141.527
$ws.Range("L134").Value = 9600
$ws.Range("M134").Value = -4606.527
$ws.Range("N134").Value = -14670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1327.3636
$ws.Range("I16").Value = 950.1667
$ws.Range("J16").Value = 1780
$ws.Range("K16").Value = 950.1667
$ws.Range("L16").Value = 1780
$ws.Range("M16").Value = -663.1667
$ws.Range("N16").Value = -2354
$ws.Range("H53").Value = 35684
$ws.Range("J53").Value = 35684
$ws.Range("L53").Value = 35684
$ws.Range("N53").Value = -36898
$ws.Range("H58").Value = 16235.5
$ws.Range("I58").Value = 1457.4445
$ws.Range("J58").Value = 32860.812
$ws.Range("K58").Value = 1457.4445
$ws.Range("L58").Value = 32860.812
$ws.Range("M58").Value = -1254.4445
$ws.Range("N58").Value = -33266.812
$ws.Range("H113").Value = 1327.3636
$ws.Range("I113").Value = 950.1667
$ws.Range("J113").Value = 1780
$ws.Range("K113").Value = 950.1667
$ws.Range("L113").Value = 1780
$ws.Range("M113").Value = 1219.8333
$ws.Range("N113").Value = -6120
$ws.Range("H115").Value = 25950
$ws.Range("J115").Value = 25950
$ws.Range("L115").Value = 25950
$ws.Range("N115").Value = -28300
$ws.Range("H122").Value = 907.16
$ws.Range("I122").Value = 887.6818
$ws.Range("J122").Value = 1050
$ws.Range("K122").Value = 2663.0454
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -213.0454
$ws.Range("N122").Value = -8050
$ws.Range("H134").Value = 1052.8889
$ws.Range("I134").Value = 938.8387
$ws.Range("K134").Value = 2816.5161
$ws.Range("M134").Value = -281.5160999999998
$ws.Range("H136").Value = 16235.5
$ws.Range("I136").Value = 1457.4445
$ws.Range("J136").Value = 32860.812
$ws.Range("K136").Value = 4372.333500000001
$ws.Range("L136").Value = 98582.43599999999
$ws.Range("M136").Value = -1822.333500000001
$ws.Range("N136").Value = -103682.436

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 83333460
$ws.Range("I38").Value = 73.333336
$ws.Range("J38").Value = 166666830
$ws.Range("K38").Value = 220.000008
$ws.Range("L38").Value = 500000490
$ws.Range("M38").Value = 126.999992
$ws.Range("N38").Value = -500001184
$ws.Range("H86").Value = 33334070
$ws.Range("I86").Value = 725.2
$ws.Range("K86").Value = 2175.6
$ws.Range("M86").Value = -989.6000000000004
$ws.Range("H89").Value = 33334070
$ws.Range("I89").Value = 725.2
$ws.Range("K89").Value = 6526.8
$ws.Range("M89").Value = -598.8000000000002
$ws.Range("H107").Value = 5735.722
$ws.Range("I107").Value = 6055.4707
$ws.Range("K107").Value = 18166.4121
$ws.Range("M107").Value = -16246.4121
$ws.Range("H113").Value = 646.55
$ws.Range("I113").Value = 508.45456
$ws.Range("J113").Value = 815.3333
$ws.Range("K113").Value = 1525.36368
$ws.Range("L113").Value = 2445.9999
$ws.Range("M113").Value = 644.6363200000001
$ws.Range("N113").Value = -6785.9999
$ws.Range("H121").Value = 1395.9231
$ws.Range("I121").Value = 577.5
$ws.Range("K121").Value = 1732.5
$ws.Range("M121").Value = -422.5
$ws.Range("H122").Value = 798.06665
$ws.Range("I122").Value = 398
$ws.Range("K122").Value = 3582
$ws.Range("M122").Value = -1132
$ws.Range("H131").Value = 104940.18
$ws.Range("J131").Value = 108296.2
$ws.Range("L131").Value = 324888.6
$ws.Range("N131").Value = -334968.6
$ws.Range("H132").Value = 411.1111
$ws.Range("I132").Value = 372.5
$ws.Range("J132").Value = 442
$ws.Range("K132").Value = 3352.5
$ws.Range("L132").Value = 3978
$ws.Range("M132").Value = -822.5
$ws.Range("N132").Value = -9038

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4818631.5
$ws.Range("I70").Value = 3666.6667
$ws.Range("J70").Value = 8945744
$ws.Range("K70").Value = 3666.6667
$ws.Range("L70").Value = 8945744
$ws.Range("M70").Value = -3396.6667
$ws.Range("N70").Value = -8946284
$ws.Range("H73").Value = 4818631.5
$ws.Range("I73").Value = 3666.6667
$ws.Range("J73").Value = 8945744
$ws.Range("K73").Value = 3666.6667
$ws.Range("L73").Value = 8945744
$ws.Range("M73").Value = -2730.6667
$ws.Range("N73").Value = -8947616
$ws.Range("H95").Value = 21007.75
$ws.Range("J95").Value = 21007.75
$ws.Range("L95").Value = 21007.75
$ws.Range("N95").Value = -26499.75
$ws.Range("H97").Value = 1537.7273
$ws.Range("I97").Value = 1207
$ws.Range("J97").Value = 2662.2
$ws.Range("K97").Value = 1207
$ws.Range("L97").Value = 2662.2
$ws.Range("M97").Value = -711
$ws.Range("N97").Value = -3654.2
$ws.Range("H122").Value = 38096164
$ws.Range("I122").Value = 13889791
$ws.Range("J122").Value = 90910060
$ws.Range("K122").Value = 41669373
$ws.Range("L122").Value = 272730180
$ws.Range("M122").Value = -41666923
$ws.Range("N122").Value = -272735080
$ws.Range("H126").Value = 5379.1724
$ws.Range("I126").Value = 4269.4736
$ws.Range("J126").Value = 7487.6
$ws.Range("K126").Value = 12808.4208
$ws.Range("L126").Value = 22462.8
$ws.Range("M126").Value = -10338.4208
$ws.Range("N126").Value = -27402.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1104.8235
$ws.Range("I22").Value = 1746.8572
$ws.Range("K22").Value = 1746.8572
$ws.Range("M22").Value = -1451.8572
$ws.Range("H27").Value = 1104.8235
$ws.Range("I27").Value = 1746.8572
$ws.Range("K27").Value = 1746.8572
$ws.Range("M27").Value = -1639.8572
$ws.Range("H46").Value = 1039.0857
$ws.Range("I46").Value = 1013.7647
$ws.Range("K46").Value = 1013.7647
$ws.Range("M46").Value = -825.7646999999999
$ws.Range("H55").Value = 172.5
$ws.Range("J55").Value = 180
$ws.Range("L55").Value = 180
$ws.Range("N55").Value = -526
$ws.Range("H111").Value = 29990.715
$ws.Range("J111").Value = 29990.715
$ws.Range("L111").Value = 29990.715
$ws.Range("N111").Value = -38170.715
$ws.Range("H135").Value = 20250
$ws.Range("J135").Value = 20250
$ws.Range("L135").Value = 20250
$ws.Range("N135").Value = -30390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 6669333.5
$ws.Range("J18").Value = 6669333.5
$ws.Range("L18").Value = 6669333.5
$ws.Range("N18").Value = -6669679.5
$ws.Range("H22").Value = 5015
$ws.Range("J22").Value = 5015
$ws.Range("L22").Value = 5015
$ws.Range("N22").Value = -5601
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H62").Value = 3285.2856
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 3499.4
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 3499.4
$ws.Range("M62").Value = -2126
$ws.Range("N62").Value = -4747.4
$ws.Range("H65").Value = 3285.2856
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 3499.4
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 17497
$ws.Range("M65").Value = -10630
$ws.Range("N65").Value = -23737

